# Applies the "Made WahWah SPI ready, added on/off feature SPI for distortion,
# updated id list of spi" edit to the workbook.
#
# Summary of the change (sheet "Blad1" / sheet1.xml):
#  - filterid 18 ("gauchissement") row: the Subsub value changes from
#    "range or off" to "range".
#  - Two new filterid rows are appended at the bottom of the table:
#      filterid 22: Subsystem = (blank, merged with "wawa" above),
#                   Subsub = "Range"
#      filterid 23: Subsystem = "gauchissement", Subsub = "on/off"
#  - The "wawa" (WahWah) Subsystem cell (row 22) is merged with the new
#    blank cell below it (row 23), just like the existing "tremololo"
#    block at rows 20:21.
#  - The active selection moves to K27 and the previous frozen/scrolled
#    topLeftCell (A6) is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. "range or off" -> "range" for the existing gauchissement row ---
$ws.Range("C19").Value = "range"

# --- 2. New row 23: filterid 22, blank Subsystem (merges with "wawa"), Subsub "Range" ---
$ws.Range("A23").Value = 22
$ws.Range("C23").Value = "Range"

# --- 3. New row 24: filterid 23, Subsystem "gauchissement", Subsub "on/off" ---
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "gauchissement"
$ws.Range("C24").Value = "on/off"

# Match the centered / vertically centered formatting used by the rest of
# the "filterid" (A) and "Subsystem" (B) columns for the new row.
$ws.Range("A24:B24").HorizontalAlignment = -4108
$ws.Range("A24:B24").VerticalAlignment = -4108

# --- 4. Merge the WahWah "Subsystem" cell down into the new row, like the
#        existing tremololo (B20:B21) block ---
$ws.Range("B22:B23").Merge() | Out-Null

# --- 5. Update the active selection / scroll position ---
$ws.Range("K27").Select() | Out-Null
